# Raul's Log - add new rows (424-442) to the "Logs" sheet, mirroring the
# rows that were logged for 9/8/2016 (serial 42621).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$eventDate = [DateTime]::FromOADate(42621)

# Column A = Task Type, B = Date, C = Time, D = Building, E = Room, F = Notes
$setupNote = "Set up and test all classrooms / test throughout   Operate 19:30-22:30"

$rows = @(
    @{ Row=424; A="Demo";       C="1600"; D="ACE"; E="003";  F=$null;      Tall=$false },
    @{ Row=425; A="Demo";       C="1600"; D="ACE"; E="011";  F=$null;      Tall=$false },
    @{ Row=426; A="Demo";       C="1600"; D="HNE"; E="032";  F=$null;      Tall=$false },
    @{ Row=427; A="Demo";       C="1600"; D="HNE"; E="401";  F=$null;      Tall=$false },
    @{ Row=428; A="Demo";       C="1900"; D="DB";  E="0004"; F=$null;      Tall=$false },
    @{ Row=429; A="Demo";       C="1900"; D="HNE"; E="031";  F=$null;      Tall=$false },
    @{ Row=430; A="Demo";       C="1900"; D="HNE"; E="401";  F=$null;      Tall=$false },
    @{ Row=431; A="Operator";   C="1930"; D="SSB"; E="E112"; F=$setupNote; Tall=$true  },
    @{ Row=432; A="Demo";       C="1930"; D="SSB"; E="E115"; F=$setupNote; Tall=$true  },
    @{ Row=433; A="Demo";       C="1930"; D="SSB"; E="E118"; F=$setupNote; Tall=$true  },
    @{ Row=434; A="Demo";       C="1930"; D="SSB"; E="N105"; F=$setupNote; Tall=$true  },
    @{ Row=435; A="Demo";       C="1930"; D="SSB"; E="N106"; F=$setupNote; Tall=$true  },
    @{ Row=436; A="Demo";       C="1930"; D="SSB"; E="N107"; F=$setupNote; Tall=$true  },
    @{ Row=437; A="Demo";       C="1930"; D="SSB"; E="N108"; F=$setupNote; Tall=$true  },
    @{ Row=438; A="Demo";       C="1930"; D="SSB"; E="W133"; F=$setupNote; Tall=$true  },
    @{ Row=439; A="Demo";       C="1930"; D="SSB"; E="W136"; F=$setupNote; Tall=$true  },
    @{ Row=440; A="Demo";       C="1930"; D="SSB"; E="E112"; F=$setupNote; Tall=$true  },
    @{ Row=441; A="Setup Mic";  C="1800"; D="DB";  E="2027"; F="Setup meck mic with small PA from DB 0003";   Tall=$false },
    @{ Row=442; A="Pickup Mic"; C="2050"; D="DB";  E="2027"; F="Return neck mic and small PA to DB 0003";     Tall=$false }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $eventDate

    if ($r.Row -eq 442) {
        # The "Return" note/time for this row was logged before its time
        # stamp was filled in, so the shared-string table picks up F
        # ("Return neck mic...") ahead of C ("2050").
        $ws.Cells.Item($rowNum, 6).Value = $r.F
        $ws.Cells.Item($rowNum, 3).Value = $r.C
        $ws.Cells.Item($rowNum, 4).Value = $r.D
        $ws.Cells.Item($rowNum, 5).Value = $r.E
    } else {
        $ws.Cells.Item($rowNum, 3).Value = $r.C
        $ws.Cells.Item($rowNum, 4).Value = $r.D
        $ws.Cells.Item($rowNum, 5).Value = $r.E

        if ($r.F) {
            $ws.Cells.Item($rowNum, 6).Value = $r.F
        }
    }

    if ($r.Tall) {
        $ws.Rows.Item($rowNum).RowHeight = 30
    }
}

# Match the position/selection left after the edits were made.
$ws.Application.Goto($ws.Range("A418"))
$ws.Range("C443").Select()
